$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "EmployeeSituation" data-string lookup table, appended below the
# existing listsCardsByClient tables (header row + two value rows).
$ws.Range("B39").Value = "ID"
$ws.Range("C39").Value = "Data String EmployeeSituation"

$ws.Range("B40").Value = "F"
$ws.Range("C40").Value = "Fijo"

$ws.Range("B41").Value = "I"
$ws.Range("C41").Value = "Identificacion"

# Give the new table the same thin box border used elsewhere in the sheet.
$ws.Range("B39:C41").Borders.LineStyle = 1
